# Update the dSF (column F) values to reflect the repulled / recalculated
# data as described by the commit message ("repull data, push all data,
# mean calculation"). Only the rows whose dSF value actually changed are
# touched; all other cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    4  = 2
    5  = 1
    6  = 9
    7  = -7
    8  = 1
    9  = 3
    11 = 3
    14 = -2
    16 = -4
    17 = -3
    18 = 2
    20 = -3
    21 = 1
    23 = -3
    25 = -7
    26 = -1
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
